# Auto-generated edit script applying the cryptos.xlsx diff
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextCell($range, [string]$value) {
    # Force text so numeric-looking / multi-dot strings keep their exact
    # characters (no coercion to a Double, no lost trailing zeros),
    # then restore the default style so we do not leave stray formatting.
    $range.NumberFormat = "@"
    $range.Value = $value
    $range.Style = "Normal"
}

# Row 2
Set-TextCell $ws.Range("D2") "62.503.11"
Set-TextCell $ws.Range("E2") "  +5.67%  "

# Row 3
Set-TextCell $ws.Range("D3") "3.463.46"
Set-TextCell $ws.Range("E3") "  +3.47%  "

# Row 4
Set-TextCell $ws.Range("E4") "  +0.02%  "

# Row 5
Set-TextCell $ws.Range("D5") "411.89"
Set-TextCell $ws.Range("E5") "  +0.23%  "

# Row 6
Set-TextCell $ws.Range("D6") "128.31"
Set-TextCell $ws.Range("E6") "  +15.01%  "

# Row 7
Set-TextCell $ws.Range("D7") "3.449.95"
Set-TextCell $ws.Range("E7") "  +3.22%  "

# Row 8
Set-TextCell $ws.Range("D8") "0.596"
Set-TextCell $ws.Range("E8") "  +1.48%  "

# Row 9
Set-TextCell $ws.Range("E9") "  +0.06%  "

# Row 10
Set-TextCell $ws.Range("D10") "0.688"

# Row 11
Set-TextCell $ws.Range("D11") "0.128"
Set-TextCell $ws.Range("E11") "  +29.50%  "

# Row 12
Set-TextCell $ws.Range("D12") "43.59"
Set-TextCell $ws.Range("E12") "  +8.50%  "

# Row 13
Set-TextCell $ws.Range("E13") "  -0.39%  "

# Row 14
Set-TextCell $ws.Range("D14") "4.000.00"
Set-TextCell $ws.Range("E14") "  +3.31%  "

# Row 15
Set-TextCell $ws.Range("D15") "8.75"
Set-TextCell $ws.Range("E15") "  +2.97%  "

# Row 16
Set-TextCell $ws.Range("D16") "20.18"
Set-TextCell $ws.Range("E16") "  +4.01%  "

# Row 17
Set-TextCell $ws.Range("D17") "3.417.57"
Set-TextCell $ws.Range("E17") "  +1.85%  "

# Row 18
Set-TextCell $ws.Range("D18") "62.433.93"
Set-TextCell $ws.Range("E18") "  +5.84%  "

# Row 19
Set-TextCell $ws.Range("E19") "  +0.00%  "

# Row 20
Set-TextCell $ws.Range("D20") "11.13"
Set-TextCell $ws.Range("E20") "  +2.36%  "

# Row 21
Set-TextCell $ws.Range("E21") "  +22.17%  "

# Row 22
Set-TextCell $ws.Range("D22") "3.36"
Set-TextCell $ws.Range("E22") "  +0.28%  "

# Row 23
Set-TextCell $ws.Range("B23") "Litecoin"
Set-TextCell $ws.Range("C23") "https://coinranking.com/coin/D7B1x_ks7WhV5+litecoin-ltc"
Set-TextCell $ws.Range("D23") "82.27"
Set-TextCell $ws.Range("E23") "  +9.06%  "

# Row 24
Set-TextCell $ws.Range("B24") "InternetComputer(DFINITY)"
Set-TextCell $ws.Range("C24") "https://coinranking.com/coin/aMNLwaUbY+internetcomputerdfinity-icp"
Set-TextCell $ws.Range("D24") "13.18"
Set-TextCell $ws.Range("E24") "  +0.36%  "

# Row 25
Set-TextCell $ws.Range("D25") "312.55"
Set-TextCell $ws.Range("E25") "  +2.68%  "

# Row 26
Set-TextCell $ws.Range("E26") "  -0.39%  "

# Row 27
Set-TextCell $ws.Range("D27") "30.37"
Set-TextCell $ws.Range("E27") "  +6.12%  "

# Row 28
Set-TextCell $ws.Range("D28") "8.20"
Set-TextCell $ws.Range("E28") "  +3.03%  "

# Row 29
Set-TextCell $ws.Range("D29") "7.82"
Set-TextCell $ws.Range("E29") "  +4.20%  "

# Row 30
Set-TextCell $ws.Range("D30") "0.121"
Set-TextCell $ws.Range("E30") "  +4.29%  "

# Row 31
Set-TextCell $ws.Range("E31") "  +3.99%  "

# Row 32
Set-TextCell $ws.Range("D32") "4.36"
Set-TextCell $ws.Range("E32") "  -2.63%  "

# Row 33
Set-TextCell $ws.Range("D33") "12.12"
Set-TextCell $ws.Range("E33") "  +3.73%  "

# Row 34
Set-TextCell $ws.Range("D34") "44.18"
Set-TextCell $ws.Range("E34") "  +9.68%  "

# Row 35
Set-TextCell $ws.Range("D35") "2.67"
Set-TextCell $ws.Range("E35") "  +25.19%  "

# Row 36
Set-TextCell $ws.Range("E36") "  +0.16%  "

# Row 37
Set-TextCell $ws.Range("D37") "0.0494"
Set-TextCell $ws.Range("E37") "  -6.32%  "

# Row 38
Set-TextCell $ws.Range("D38") "52.67"
Set-TextCell $ws.Range("E38") "  +1.04%  "

# Row 39
Set-TextCell $ws.Range("E39") "  +2.59%  "

# Row 40
Set-TextCell $ws.Range("D40") "0.996"
Set-TextCell $ws.Range("E40") "  -0.34%  "

# Row 41
Set-TextCell $ws.Range("E41") "  -3.35%  "

# Row 42
Set-TextCell $ws.Range("E42") "  +3.63%  "

# Row 43
Set-TextCell $ws.Range("E43") "  +2.76%  "

# Row 44
Set-TextCell $ws.Range("D44") "137.85"
Set-TextCell $ws.Range("E44") "  -0.78%  "

# Row 45
Set-TextCell $ws.Range("D45") "17.87"
Set-TextCell $ws.Range("E45") "  +5.38%  "

# Row 46
Set-TextCell $ws.Range("B46") "NEARProtocol"
Set-TextCell $ws.Range("C46") "https://coinranking.com/coin/DCrsaMv68+nearprotocol-near"
Set-TextCell $ws.Range("D46") "4.03"
Set-TextCell $ws.Range("E46") "  +0.35%  "

# Row 47
Set-TextCell $ws.Range("B47") "TheGraph"
Set-TextCell $ws.Range("C47") "https://coinranking.com/coin/qhd1biQ7M+thegraph-grt"
Set-TextCell $ws.Range("D47") "0.290"
Set-TextCell $ws.Range("E47") "  +3.58%  "

# Row 48
Set-TextCell $ws.Range("E48") "  +0.87%  "

# Row 49
Set-TextCell $ws.Range("D49") "22.46"
Set-TextCell $ws.Range("E49") "  +0.01%  "

# Row 50
Set-TextCell $ws.Range("D50") "2.223.46"
Set-TextCell $ws.Range("E50") "  +0.81%  "

# Row 51
Set-TextCell $ws.Range("D51") "3.799.57"
Set-TextCell $ws.Range("E51") "  +3.27%  "
